$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at row 51 (old rows 51-62 shift down to 52-63) ---
# new entry: "2005 - 京都市衛生公害研究所"
$ws.Rows.Item(51).Insert()

$a51 = $ws.Range("A51")
$a51.NumberFormat = "@"       # force text storage so "2005" isn't read back as a number
$a51.Value = '2005'
$a51.ClearFormats()           # drop the temporary text format, back to default styling

$ws.Range("B51").Value = '**京都市衛生公害研究所** <br> [市販ナチュラルチーズからのリステリア菌の検出](https://www.city.kyoto.lg.jp/hokenfukushi/cmsfiles/contents/0000118/118277/o-5.pdf)'
$ws.Range("C51").Value = '済'

# --- Insert a second new row at row 58 (after old row 56, now shifted to row 57) ---
# shifts old rows 57-62 (now at 58-63) further down to 59-64
# new entry: "2001 - 神戸市環境保健研究所"
$ws.Rows.Item(58).Insert()

$a58 = $ws.Range("A58")
$a58.NumberFormat = "@"
$a58.Value = '2001'
$a58.ClearFormats()

$ws.Range("B58").Value = '**神戸市環境保健研究所** <br> [市販の輸入生野菜および果物における病原菌汚染の実態調査](https://www.jstage.jst.go.jp/article/jsfm1994/19/2/19_2_71/_pdf/-char/ja) <br>（日本食品微生物学会雑誌, 19(2), 71-75, 2002）'
$ws.Range("C58").Value = '未登録'
